$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap country labels for rows 41 and 42 (Republica Dominicana overtook Panama in ranking)
$ws.Range("A41").Value = "Republica Dominicana"
$ws.Range("A42").Value = "Panama"

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Julio de 2020 a las 17:44"

# Row 4 (rank 8) - Rusia
$ws.Range("B4").Value = 3978665
$ws.Range("C4").Value = 17236
$ws.Range("D4").Value = 1851157
$ws.Range("E4").Value = 1983313
$ws.Range("G4").Value = 361
$ws.Range("H4").Value = 144195

# Row 5 (rank 9) - Sudafrica
$ws.Range("B5").Value = 2122577
$ws.Range("C5").Value = 932
$ws.Range("E5").Value = 633111
$ws.Range("G5").Value = 13
$ws.Range("H5").Value = 80264

# Row 6 (rank 10) - Peru
$ws.Range("B6").Value = 1171446
$ws.Range("C6").Value = 16529
$ws.Range("D6").Value = 737808
$ws.Range("E6").Value = 405150
$ws.Range("G6").Value = 389
$ws.Range("H6").Value = 28488

# Row 13 (rank 17) - Turquia
$ws.Range("B13").Value = 295817
$ws.Range("C13").Value = 445
$ws.Range("G13").Value = 110
$ws.Range("H13").Value = 45422

# Row 17 (rank 21) - Canada
$ws.Range("B17").Value = 244752
$ws.Range("C17").Value = 128
$ws.Range("D17").Value = 197431
$ws.Range("E17").Value = 12248
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 35073

# Row 21 (rank 25)
$ws.Range("D21").Value = 188100
$ws.Range("E21").Value = 6282

# Row 41 (rank 45) - now Republica Dominicana
$ws.Range("B41").Value = 54797
$ws.Range("C41").Value = 841
$ws.Range("D41").Value = 25976
$ws.Range("E41").Value = 27822
$ws.Range("G41").Value = 6
$ws.Range("H41").Value = 999

# Row 42 (rank 46) - now Panama
$ws.Range("B42").Value = 54426
$ws.Range("D42").Value = 29164
$ws.Range("E42").Value = 24135
$ws.Range("H42").Value = 1127

# Row 46 (rank 50)
$ws.Range("D46").Value = 44584
$ws.Range("E46").Value = 3823

# Row 70 (rank 74)
$ws.Range("B70").Value = 14160
$ws.Range("C70").Value = 62
$ws.Range("D70").Value = 8899
$ws.Range("E70").Value = 4901
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 360

# Row 86 (rank 90)
$ws.Range("B86").Value = 8787
$ws.Range("C86").Value = 308
$ws.Range("E86").Value = 4520

# Row 97 (rank 101)
$ws.Range("B97").Value = 5027
$ws.Range("C97").Value = 7
$ws.Range("D97").Value = 4884
$ws.Range("E97").Value = 87

# Row 116 (rank 120)
$ws.Range("B116").Value = 2449
$ws.Range("C116").Value = 3
$ws.Range("D116").Value = 2321
$ws.Range("E116").Value = 41

# Row 135 (rank 139)
$ws.Range("B135").Value = 1536
$ws.Range("C135").Value = 29
$ws.Range("E135").Value = 1020
